$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: new row with M2 Standoffs name (create this shared string first so it
# lands at index 57, matching the order new strings were appended upstream).
$ws.Range("B20").Value = "M2 Standoffs "

# Row 19 already has B19 = "M2 Heat-Set Knurled Threaded Inserts"; add missing hyperlink in G19.
# Pre-apply the built-in "Hyperlink" cell style so the style the engine bakes
# in for the new hyperlink cell reuses a single, compact cellXfs entry.
$ws.Range("G19").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G19"), "https://www.amazon.com/AIEX-Printing-Embedment-Automotive-M2x3x3-5mm/dp/B0B8GN63S2/ref=sr_1_3?keywords=m2+threaded+insert&qid=1667426629&qu=eyJxc2MiOiI0LjUwIiwicXNhIjoiMy44NSIsInFzcCI6IjMuNTYifQ%3D%3D&sprefix=m2+threaded%2Caps%2C147&sr=8-3")

# Row 20: add vendor link for M2 Standoffs
$ws.Range("G20").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G20"), "https://www.amazon.com/HELIFOUNER-Spacers-Standoffs-Assortment-Tweezers/dp/B09F8TCLRY/ref=sr_1_4?keywords=m2%2Bstandoffs&qid=1667426997&qu=eyJxc2MiOiIzLjg1IiwicXNhIjoiMy4xNSIsInFzcCI6IjMuMDYifQ%3D%3D&s=industrial&sr=1-4&th=1")

# Match the author's final selection/active cell (G20, the last-edited cell).
[void]$ws.Range("G20").Select()
